$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing row 1 values (total APF area update)
$ws.Range("A1").Value = 300
$ws.Range("C1").Value = 50
$ws.Range("D1").Value = 50

# Append new rows 4-7
$ws.Range("A4").Value = 300
$ws.Range("B4").Value = 300
$ws.Range("C4").Value = 50
$ws.Range("D4").Value = 50
$ws.Range("E4").Value = 0

$ws.Range("A5").Value = 300
$ws.Range("B5").Value = 350
$ws.Range("C5").Value = 50
$ws.Range("D5").Value = 50
$ws.Range("E5").Value = 0

$ws.Range("A6").Value = 300
$ws.Range("B6").Value = 200
$ws.Range("C6").Value = 50
$ws.Range("D6").Value = 50
$ws.Range("E6").Value = 0

$ws.Range("A7").Value = 300
$ws.Range("B7").Value = 150
$ws.Range("C7").Value = 50
$ws.Range("D7").Value = 50
$ws.Range("E7").Value = 0

# Update selection to match the final committed sheet view
# (user extended the selection from A1 down through the newly-added rows to E7)
$ws.Range("A1:E7").Select()
